$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -11.931
$ws.Range("B7").Value = 5.515
$ws.Range("A8").Value = -22.03
$ws.Range("A10").Value = -21.223
$ws.Range("D10").Value = -8.068999999999999
$ws.Range("A12").Value = -21.593
$ws.Range("D12").Value = -7.269
$ws.Range("D13").Value = -8.715
$ws.Range("D14").Value = -7.85
$ws.Range("B15").Value = 5.236
$ws.Range("E16").Value = 16.764
$ws.Range("A18").Value = -21.105
$ws.Range("B18").Value = 8.114999999999998
$ws.Range("C18").Value = -11.835
$ws.Range("E18").Value = 17.386
$ws.Range("C19").Value = -12.051
$ws.Range("B20").Value = 6.956999999999999
$ws.Range("E22").Value = 16.581
$ws.Range("E26").Value = 16.977
$ws.Range("C27").Value = -13.407
$ws.Range("B29").Value = 5.749000000000001
$ws.Range("D29").Value = -7.633
$ws.Range("B30").Value = 5.365
$ws.Range("B31").Value = 5.167
$ws.Range("C31").Value = -12.766
$ws.Range("D32").Value = -8.274000000000001
$ws.Range("D35").Value = -7.386999999999999
$ws.Range("A37").Value = -20.02
$ws.Range("C38").Value = -12.845
$ws.Range("E39").Value = 16.583
$ws.Range("B40").Value = 8.901
$ws.Range("C42").Value = -12.328
$ws.Range("D43").Value = -8.471999999999998
$ws.Range("C44").Value = -12.741
$ws.Range("E44").Value = 16.756
$ws.Range("C47").Value = -12.335
$ws.Range("D48").Value = -7.793000000000001
$ws.Range("D49").Value = -8.071000000000002
$ws.Range("B50").Value = 5.867
$ws.Range("D50").Value = -8.085999999999999
$ws.Range("E51").Value = 16.665
$ws.Range("E54").Value = 16.551
$ws.Range("A55").Value = -21.868
$ws.Range("D56").Value = -8.434000000000001
$ws.Range("E57").Value = 16.45
$ws.Range("C58").Value = -12.727
$ws.Range("E63").Value = 17.659
$ws.Range("C65").Value = -12.567
$ws.Range("A68").Value = -21.681
$ws.Range("B68").Value = 5.348000000000001
$ws.Range("D69").Value = -7.766999999999999
$ws.Range("C73").Value = -12.77
$ws.Range("B76").Value = 6.556999999999999
$ws.Range("A77").Value = -19.876
$ws.Range("E77").Value = 16.918
$ws.Range("A78").Value = -19.774
$ws.Range("A81").Value = -21.709
$ws.Range("D81").Value = -7.877000000000001
$ws.Range("A82").Value = -21.896
$ws.Range("E86").Value = 16.431
$ws.Range("B87").Value = 4.88
$ws.Range("B88").Value = 5.891
$ws.Range("C90").Value = -13.331
$ws.Range("D92").Value = -7.175
$ws.Range("C94").Value = -11.085
$ws.Range("C95").Value = -11.928
$ws.Range("B96").Value = 6.581999999999999
$ws.Range("E96").Value = 16.541
$ws.Range("B98").Value = 5.736
$ws.Range("E98").Value = 16.455
$ws.Range("B101").Value = 7.716999999999999
$ws.Range("C101").Value = -12.978
$ws.Range("B102").Value = 8.009
